# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet and the per-locale "Status" column.
#  - Latest Handback DateTime is refreshed for zh-cn and de-de.
#  - The stale "version mismatch" Error Detail messages are cleared now that
#    the handback is in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus
$ws.Range("E3").Value = $newStatus
$ws.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de status columns so the longer text fits.
$ws.Columns.Item(5).ColumnWidth = 29.1
$ws.Columns.Item(6).ColumnWidth = 29.1

# --- zh-cn sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = $newStatus
$ws.Range("C3").Value = $newStatus

$ws.Range("K2").Value = "2016-09-07 13:14:19"
$ws.Range("K3").Value = "2016-09-07 13:14:19"

$ws.Range("P2").Value = ""
$ws.Range("P3").Value = ""

$ws.Columns.Item(3).ColumnWidth = 29.1
$ws.Columns.Item(16).ColumnWidth = 12.75

# --- de-de sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = $newStatus
$ws.Range("C3").Value = $newStatus

$ws.Range("K2").Value = "2016-09-07 13:14:37"
$ws.Range("K3").Value = "2016-09-07 13:14:37"

$ws.Range("P2").Value = ""
$ws.Range("P3").Value = ""

$ws.Columns.Item(3).ColumnWidth = 29.1
$ws.Columns.Item(16).ColumnWidth = 12.75
